$d = $word.ActiveDocument

# --- Paragraph 1: "Nuget" + "-package manager:" -> single run "Nuget-package manager:"
# (also removes the now-orphaned proofErr spellStart/spellEnd pair around "Nuget")
# A throw-away marker character is inserted at the very start of the document so the
# Find/Replace match no longer begins exactly at the paragraph's first character; that
# makes the engine fold the split runs together and drop the bracketing <w:proofErr/>.
$r0 = $d.Range(0, 0)
$r0.InsertBefore("Z")
$d.Content.Find.Execute("ZNuget-package manager:", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Nuget-package manager:", 2)

# Replacing across the run boundary drops the run's character formatting, so restore it
# to match the paragraph's original Consolas/9.5pt/black look.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Font.Name = "Consolas"
$p1.Range.Font.NameAscii = "Consolas"
$p1.Range.Font.NameBi = "Consolas"
$p1.Range.Font.Size = 9.5
$p1.Range.Font.SizeBi = 9.5
$p1.Range.Font.Color = 0

# --- Paragraph 2: "install-package " + "EntityFramework" -> single run
# "install-package EntityFramework"
# (also removes the now-orphaned proofErr spellStart/spellEnd pair around "EntityFramework")
# Same trick, but the marker goes right after the match so the boundary shifts away from
# the trailing proofErr tag instead.
$p2 = $d.Paragraphs.Item(2)
$pEnd = $p2.Range.End - 1
$rIns = $d.Range($pEnd, $pEnd)
$rIns.InsertAfter("Z")
$d.Content.Find.Execute("install-package EntityFrameworkZ", $true, $false, $false, $false, `
    $false, $true, 1, $false, "install-package EntityFramework", 2)

# --- Add four new Consolas command-line paragraphs after "install-package EntityFramework"
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = "enable-migrations"

$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Text = "add-migration name"

$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Text = "add-migration name -force // to override"

$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Text = "update-database"
